$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ausgaben")

$ws.Range("A16").Value = "Lohn"
$ws.Range("D16").Formula = "=1500*1.1"
$ws.Range("D17:D25").Formula = "=1500*1.1"

$ws.Range("D27").Select()
